# test data updation for skip batch2
# Updates the pCloudy endpoint on the Capabilities sheet and refreshes the
# device roster on the DeviceList sheet (6 devices trimmed down to 5, with
# new device/version/credential data), matching the "SkippedRerunBatches2"
# -> "DBS_Automation/10" batch re-run.

$wb = $excel.ActiveWorkbook

# ---- Capabilities sheet -----------------------------------------------
$caps = $wb.Worksheets.Item("Capabilities")
$caps.Range("D2").Value = "'https://ind-west.pcloudy.com"
$caps.Range("D3").Value = "'https://ind-west.pcloudy.com"
$caps.Application.ActiveWindow.RangeSelection
$caps.Range("D7").Select()

# ---- DeviceList sheet ---------------------------------------------------
$dl = $wb.Worksheets.Item("DeviceList")

# Drop the last six device columns (G:L) - only 5 devices remain (B:F)
$dl.Range("G1:L10").Delete(-4159)

# New device names (row 1)
$dl.Range("B1").Value = "XIAOMI_RedmiGo_Android_8.1.0_1d174"
$dl.Range("C1").Value = "REALME_8s_Android_11.0.0_403e0"
$dl.Range("D1").Value = "ONEPLUS_Nord2_Android_11.0.0_9e5b2"
$dl.Range("E1").Value = "XIAOMI_RedmiNote8_Android_11.0.0_da311"
$dl.Range("F1").Value = "VIVO_Y15_Android_9.0.0_6bc8e"

# New device OS versions (row 2)
$dl.Range("B2").Value = "8.1.0"
$dl.Range("C2").Value = "11.0.0"
$dl.Range("D2").Value = "11.0.0"
$dl.Range("F2").Value = "9.0.0"

# New UserName credentials (row 4)
$dl.Range("B4").Value = "S2325475AUID"
$dl.Range("C4").Value = "S2325476ZUID"
$dl.Range("D4").Value = "S2325477HUID"
$dl.Range("E4").Value = "S2325481FUID"
$dl.Range("F4").Value = "S2325482DUID"

$dl.Range("F15").Select()
